$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "вул. Горького, 22, " -> "49000, м. Дніпро"
# (formatting of the run is unchanged: bCs, sz22, szCs22)
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(9)
$full1 = $p1.Range
$r1 = $d.Range($full1.Start, $full1.End - 1)
$r1.Text = "49000, м. Дніпро"

# ---------------------------------------------------------------------------
# Change 2: paragraph "49000, м. Дніпро" -> "вул. Княгині Ольги, 22"
# formatting changes: drop bCs, add color 222222 + shd clear/auto/FFFFFF
# on both the paragraph mark run properties and the run itself.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(10)
$full2 = $p2.Range
$xml2 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:pPr><w:ind w:left='5760'/><w:rPr>" +
        "<w:color w:val='222222'/><w:sz w:val='22'/><w:szCs w:val='22'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
        "</w:rPr></w:pPr>" +
        "<w:r><w:rPr>" +
        "<w:color w:val='222222'/><w:sz w:val='22'/><w:szCs w:val='22'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
        "</w:rPr><w:t>вул. Княгині Ольги, 22</w:t></w:r></w:p>"
$full2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 3: split the "16 серпня 2018" date into a {5} placeholder run,
# keeping the same run formatting (Times New Roman / bCs / lang uk-UA)
# for all three resulting runs.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(18)
$full3 = $p3.Range
$r3 = $d.Range($full3.Start, $full3.End - 1)
$xml3 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:r><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman'/><w:bCs/><w:lang w:val='uk-UA'/></w:rPr>" +
        "<w:t xml:space='preserve'>Вказане рішення було подано державному реєстратору </w:t></w:r>" +
        "<w:r><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman'/><w:bCs/><w:lang w:val='uk-UA'/></w:rPr>" +
        "<w:t>{5}</w:t></w:r>" +
        "<w:r><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman'/><w:bCs/><w:lang w:val='uk-UA'/></w:rPr>" +
        "<w:t xml:space='preserve'> року для проведення реєстраційної дії «Внесення рішення засновників (учасників) юридичної особи або уповноваженого ними органу щодо припинення юридичної особи». </w:t></w:r>" +
        "</w:p>"
$r3.InsertXML($xml3)

Write-Output "done"
